$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'241.35"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'21.81"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.339"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'3.418"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'6.287"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.8037"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.8450"
$c.Style = "Normal"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Range("D10")
$c.Value = "'0.01076"
$c.Style = "Normal"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D11")
$c.Value = "'0.1432"
$c.Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D12")
$c.Value = "'0.07265"
$c.Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c = $ws.Range("D13")
$c.Value = "'0.03052"
$c.Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D14")
$c.Value = "'0.03148"
$c.Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D15")
$c.Value = "'0.09347"
$c.Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$c = $ws.Range("D16")
$c.Value = "'3.909"
$c.Style = "Normal"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D17")
$c.Value = "'0.001586"
$c.Style = "Normal"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$c = $ws.Range("D18")
$c.Value = "'0.04800"
$c.Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"
$c = $ws.Range("D19")
$c.Value = "'0.006346"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'0.0009998"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'0.004050"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'2.164"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.3233"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'0.0003031"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.03815"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.006744"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.1046"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.003204"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.006206"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.00005609"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.5807"
$c.Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$c = $ws.Range("D48")
$c.Value = "'0.1427"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.01011"
$c.Style = "Normal"
